# Append daily COVID boletim rows for 2021-02-25 through 2021-03-02
# (serial dates 44222-44257) to Planilha1, columns A (date) through J.
$rows = @(
    @(44222,5029,119,1096,6244,984,89,5,84,23),
    @(44223,5039,130,1103,6272,994,86,5,81,23),
    @(44224,5097,141,1130,6368,1004,103,5,98,23),
    @(44225,5145,112,1142,6399,1012,107,5,102,23),
    @(44226,5172,95,1146,6413,1019,103,6,97,24),
    @(44227,5189,83,1147,6419,1030,92,6,86,25),
    @(44228,5203,109,1149,6461,1044,79,7,72,26),
    @(44229,5212,144,1154,6510,1053,75,7,68,26),
    @(44230,5219,200,1160,6579,1066,68,8,60,26),
    @(44231,5235,206,1163,6604,1076,60,9,51,27),
    @(44232,5313,130,1196,6639,1090,79,12,67,27),
    @(44233,5322,118,1200,6640,1096,77,12,65,27),
    @(44234,5325,120,1204,6649,1099,77,11,66,28),
    @(44235,5380,98,1209,6687,1112,69,11,58,28),
    @(44236,5397,116,1210,6723,1118,64,12,52,28),
    @(44237,5438,104,1213,6755,1126,59,11,48,28),
    @(44238,5482,114,1225,6821,1128,68,10,58,29),
    @(44239,5501,126,1227,6854,1135,63,10,53,29),
    @(44240,5529,113,1237,6879,1148,59,11,48,30),
    @(44241,5531,117,1237,6885,1151,55,10,45,31),
    @(44242,5555,88,1245,6888,1157,57,8,49,31),
    @(44243,5557,108,1247,6912,1160,56,7,49,31),
    @(44244,5575,110,1255,6940,1166,58,8,50,31),
    @(44245,5592,120,1262,6974,1167,64,8,56,31),
    @(44246,5608,146,1267,7021,1176,60,8,52,31),
    @(44247,5633,158,1281,7072,1178,72,8,64,31),
    @(44248,5682,88,1303,7073,1180,92,8,84,31),
    @(44249,5695,88,1304,7087,1183,90,7,83,31),
    @(44250,5704,127,1325,7156,1183,110,8,102,32),
    @(44251,5736,226,1333,7295,1197,103,7,96,33),
    @(44252,5760,207,1362,7329,1203,126,7,119,33),
    @(44253,5808,249,1372,7429,1219,120,9,111,33),
    @(44254,5841,294,1388,7523,1228,127,7,120,33),
    @(44255,5847,302,1389,7538,1228,128,8,120,33),
    @(44256,5972,186,1428,7586,1233,162,8,154,33),
    @(44257,6012,218,1433,7663,1233,167,7,160,33)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$startRow = 209
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Update the active selection to reflect where the user left off editing
$ws.Range("B245").Select()

Write-Host "Added $($rows.Count) rows (209-244) to Planilha1"
